# Update sheet title / caption for the new "through" date (10-19 -> 10-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2021-10-20"

# Update the label in A11 (October row)
$ws.Range("A11").Value = "October (through 10-20)"

# H8 (July 2021) count revised down by 1
$ws.Range("H8").Value = 149

# Row 11 - October monthly counts updated
$ws.Range("B11").Value = 19
$ws.Range("C11").Value = 32
$ws.Range("D11").Value = 35
$ws.Range("E11").Value = 48
$ws.Range("F11").Value = 31
$ws.Range("G11").Value = 93
$ws.Range("H11").Value = 127

# Row 12 - Total row updated to reflect the new October figures
$ws.Range("B12").Value = 245
$ws.Range("C12").Value = 461
$ws.Range("D12").Value = 662
$ws.Range("E12").Value = 596
$ws.Range("F12").Value = 453
$ws.Range("G12").Value = 994
$ws.Range("H12").Value = 1374
